$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Append the new day (2020-05-12, serial 43963) to the three data sheets.
#    A new row is inserted just above the trailing footer/note row on each
#    sheet, which shifts that footer row down by one and picks up the
#    number formats / styles from the row above automatically.
# ---------------------------------------------------------------------------

# "all" sheet -> new row 35 (footer moves 35 -> 36)
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(35).Insert()
$wsAll.Range("A35").Value = 43963
$wsAll.Range("B35").Value = 278
$wsAll.Range("C35").Value = 275
$wsAll.Range("D35").Value = 80
$wsAll.Range("E35").Value = 69
$wsAll.Range("F35").Value = 11
$wsAll.Range("G35").Value = 9
$wsAll.Range("H35").Value = 186

# "kobe" sheet -> new row 90 (footer moves 90 -> 91)
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows.Item(90).Insert()
$wsKobe.Range("A90").Value = 43963
$wsKobe.Range("B90").Value = 0
$wsKobe.Range("C90").Value = 2659
$wsKobe.Range("D90").Value = 0
$wsKobe.Range("E90").Value = 278
$wsKobe.Range("F90").Value = 75
$wsKobe.Range("G90").Value = 65
$wsKobe.Range("H90").Value = 10
$wsKobe.Range("I90").Value = 9
$wsKobe.Range("J90").Value = 177

# "other" sheet -> new row 65 (footer moves 65 -> 66)
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(65).Insert()
$wsOther.Range("A65").Value = 43963
$wsOther.Range("B65").Value = 0
$wsOther.Range("C65").Value = 14
$wsOther.Range("D65").Value = 5
$wsOther.Range("E65").Value = 4
$wsOther.Range("F65").Value = 1
$wsOther.Range("G65").Value = 0
$wsOther.Range("H65").Value = 9

# ---------------------------------------------------------------------------
# 2) The bot's paste pulled in a stray external-workbook reference (an
#    orphaned link to a "検査件数 (累計)" sheet with no live formula left
#    behind). Reproduce that by briefly writing an external-reference
#    formula on a scratch sheet, then discard the scratch sheet - the
#    external link definition itself survives in the saved workbook.
# ---------------------------------------------------------------------------
$scratch = $wb.Worksheets.Add()
$scratch.Range("A1").Formula = "='[book2.xlsx]検査件数 (累計)'!A1"
$null = $scratch.Delete()

# ---------------------------------------------------------------------------
# 3) View-state updates: the active tab moves from "kobe" to "all", and the
#    selected cell on each touched sheet moves to the new last data row.
# ---------------------------------------------------------------------------
$wsOther.Activate()
$wsOther.Range("I65").Select()

$wsKobe.Activate()
$wsKobe.Range("J90").Select()

$wsAll.Activate()
$wsAll.Range("I35").Select()
